$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 98, shifting rows 98-210 down to 99-211.
$ws.Rows.Item(98).Insert()

# Populate the new row 98 with the new data record.
$ws.Cells.Item(98, 1).Value = 10
$ws.Cells.Item(98, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(98, 3).Value = "La Araucanía"
$ws.Cells.Item(98, 4).Value = 44601
$ws.Cells.Item(98, 5).Value = 9
$ws.Cells.Item(98, 6).Value = 100112039
$ws.Cells.Item(98, 7).Value = "Ciboulette"
$ws.Cells.Item(98, 8).Value = "Sin especificar"
$ws.Cells.Item(98, 9).Value = "Primera"
$ws.Cells.Item(98, 10).Value = 65
$ws.Cells.Item(98, 11).Value = 5000
$ws.Cells.Item(98, 12).Value = 5000
$ws.Cells.Item(98, 13).Value = 5000
$ws.Cells.Item(98, 14).Value = "$/docena de atados"
$ws.Cells.Item(98, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(98, 16).Value = 1667
$ws.Cells.Item(98, 17).Value = 3
$ws.Cells.Item(98, 18).Value = "Hortaliza"

# Apply the same date number format used by all other date cells in column D
$ws.Cells.Item(98, 4).NumberFormat = $ws.Cells.Item(99, 4).NumberFormat()
